$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Depositos")
$ws.Activate()
$ws.Range("N4").Select()
